# Updates SAFEkits.xlsx:
#  - Fill in previously-blank 2025 monthly figures on "Resumen"
#  - Update several 2025 totals across sheets
#  - Add a new "Región Policiaca" sheet with regional breakdown data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Resumen sheet: fill in missing 2025 figures (Kit sin querella count,
#    column E) and correct a couple of Kit con querella values.
# ---------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

$resumen.Range("E41").Value = 14
$resumen.Range("E42").Value = 20
$resumen.Range("E43").Value = 15
$resumen.Range("C44").Value = 18
$resumen.Range("E44").Value = 10
$resumen.Range("K44").Value = 1
$resumen.Range("E45").Value = 13
$resumen.Range("M41").Value = 3
$resumen.Range("M42").Value = 1
$resumen.Range("M43").Value = 5
$resumen.Range("M44").Value = 1
$resumen.Range("M45").Value = 0
$resumen.Range("C46").Value = 18
$resumen.Range("E46").Value = 16
$resumen.Range("K46").Value = 1
$resumen.Range("M46").Value = 1

$resumen.Range("B55").Value = 190
$resumen.Range("J55").Value = 20

# ---------------------------------------------------------------------
# 2) Data sheet: update 2025 totals
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")
$data.Range("B8").Value = 190
$data.Range("C8").Value = 20

# ---------------------------------------------------------------------
# 3) Edades sheet: update 2025 totals
# ---------------------------------------------------------------------
$edades = $wb.Worksheets.Item("Edades")
$edades.Range("B8").Value = 89
$edades.Range("C8").Value = 101

# ---------------------------------------------------------------------
# 4) Kit Analizados sheet: update 2025 totals
# ---------------------------------------------------------------------
$kitAnalizados = $wb.Worksheets.Item("Kit Analizados")
$kitAnalizados.Range("B6").Value = 98
$kitAnalizados.Range("D6").Value = 181

# ---------------------------------------------------------------------
# 5) New sheet: Región Policiaca
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$regionSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$regionSheet.Name = "Región Policiaca"

$regionSheet.Columns.Item(1).ColumnWidth = 16

$regionSheet.Range("A1").Value = "Región Policiaca"
$regionSheet.Range("B1").Value = 2021
$regionSheet.Range("C1").Value = 2022
$regionSheet.Range("D1").Value = 2023
$regionSheet.Range("E1").Value = 2024
$regionSheet.Range("F1").Value = 2025

$regiones = @(
    @("Aguadilla", 17, 14, 10, 13, 9),
    @("Aibonito", 5, 6, 4, 3, 4),
    @("Arecibo", 8, 9, 10, 17, 11),
    @("Bayamón", 20, 27, 24, 30, 39),
    @("Caguas", 12, 20, 11, 23, 17),
    @("Carolina", 8, 15, 18, 14, 17),
    @("Fajardo", 4, 10, 3, 9, 9),
    @("Guayama", 3, 2, 1, 3, 5),
    @("Humacao", 5, 7, 2, 8, 6),
    @("Mayagüez", 11, 13, 8, 10, 6),
    @("Ponce", 10, 12, 23, 9, 19),
    @("San Juan", 35, 34, 29, 45, 44),
    @("Utuado", 3, 1, 2, 5, 4)
)

$row = 2
foreach ($r in $regiones) {
    $regionSheet.Cells.Item($row, 1).Value = $r[0]
    $regionSheet.Cells.Item($row, 2).Value = $r[1]
    $regionSheet.Cells.Item($row, 3).Value = $r[2]
    $regionSheet.Cells.Item($row, 4).Value = $r[3]
    $regionSheet.Cells.Item($row, 5).Value = $r[4]
    $regionSheet.Cells.Item($row, 6).Value = $r[5]
    $row++
}

# Total row
$regionSheet.Range("A15").Value = "Total"
$regionSheet.Range("B15").Formula = "=SUM(B2:B14)"
$regionSheet.Range("C15").Formula = "=SUM(C2:C14)"
$regionSheet.Range("D15").Formula = "=SUM(D2:D14)"
$regionSheet.Range("E15").Formula = "=SUM(E2:E14)"
$regionSheet.Range("F15").Formula = "=SUM(F2:F14)"

$regionSheet.Range("A15:F15").Font.Bold = $true

$regionSheet.Activate()
$regionSheet.Range("B14").Select()
